$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "292.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.42%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.57%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.974"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.61%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07464"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.87%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.244"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-4.41%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.771"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9185"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.48%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09416"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "19.00%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1740"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.77%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08259"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.81%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03274"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.63%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09925"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.12%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001499"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.24%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005724"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.11%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.467"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.13%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.775"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.80%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.195"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "5.83%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3331"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.34%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.089"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.58%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.17%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04535"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.32%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001219"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004308"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.38%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001299"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.08%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003392"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.08%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01617"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.02%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04583"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.49%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007463"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.86%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009834"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.54%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1358"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.59%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002219"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.69%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008934"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.43%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006103"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.08%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.654"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "18.44%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001999"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-31.06%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.09%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.09%"
